$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 209.12
$ws.Range("I15").Value = 209.12
$ws.Range("K15").Value = 627.36
$ws.Range("M15").Value = -458.36

$ws.Range("H45").Value = 500
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 1500
$ws.Range("M45").Value = -1308

$ws.Range("H74").Value = 5660720.5
$ws.Range("I74").Value = 6367885.5
$ws.Range("K74").Value = 6367885.5
$ws.Range("M74").Value = -6366949.5

$ws.Range("H77").Value = 5660720.5
$ws.Range("I77").Value = 6367885.5
$ws.Range("K77").Value = 31839427.5
$ws.Range("M77").Value = -31834747.5

$ws.Range("H129").Value = 974.4722
$ws.Range("I129").Value = 391.375
$ws.Range("J129").Value = 1141.0714
$ws.Range("K129").Value = 1174.125
$ws.Range("L129").Value = 3423.2142
$ws.Range("M129").Value = 3825.875
$ws.Range("N129").Value = -13423.2142

$ws.Range("H137").Value = 802.46155
$ws.Range("I137").Value = 756.875
$ws.Range("J137").Value = 1349.5
$ws.Range("K137").Value = 2270.625
$ws.Range("L137").Value = 4048.5
$ws.Range("M137").Value = 279.375
$ws.Range("N137").Value = -9148.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 233
$ws.Range("I5").Value = 79.59999999999999
$ws.Range("K5").Value = 79.59999999999999
$ws.Range("M5").Value = 32.40000000000001

$ws.Range("H32").Value = 2777.18
$ws.Range("I32").Value = 2777.18
$ws.Range("K32").Value = 2777.18
$ws.Range("M32").Value = -2490.18

$ws.Range("H61").Value = 1523.8334
$ws.Range("I61").Value = 1501.3125
$ws.Range("J61").Value = 1704
$ws.Range("K61").Value = 1501.3125
$ws.Range("L61").Value = 1704
$ws.Range("M61").Value = -1289.3125
$ws.Range("N61").Value = -2128

$ws.Range("H74").Value = 989.4474
$ws.Range("I74").Value = 967.6177
$ws.Range("J74").Value = 1175
$ws.Range("K74").Value = 967.6177
$ws.Range("L74").Value = 1175
$ws.Range("M74").Value = -93.61770000000001
$ws.Range("N74").Value = -2923

$ws.Range("H77").Value = 989.4474
$ws.Range("I77").Value = 967.6177
$ws.Range("J77").Value = 1175
$ws.Range("K77").Value = 4838.0885
$ws.Range("L77").Value = 5875
$ws.Range("M77").Value = -470.0884999999998
$ws.Range("N77").Value = -14611

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H122").Value = 1334.0714
$ws.Range("I122").Value = 1173
$ws.Range("J122").Value = 1398.5
$ws.Range("K122").Value = 3519
$ws.Range("L122").Value = 4195.5
$ws.Range("M122").Value = -1069
$ws.Range("N122").Value = -9095.5

$ws.Range("H132").Value = 1324.3158
$ws.Range("I132").Value = 1237.5555
$ws.Range("J132").Value = 1537.2727
$ws.Range("K132").Value = 3712.6665
$ws.Range("L132").Value = 4611.8181
$ws.Range("M132").Value = -1182.6665
$ws.Range("N132").Value = -9671.8181

$ws.Range("H136").Value = 1523.8334
$ws.Range("I136").Value = 1501.3125
$ws.Range("J136").Value = 1704
$ws.Range("K136").Value = 4503.9375
$ws.Range("L136").Value = 5112
$ws.Range("M136").Value = -1953.9375
$ws.Range("N136").Value = -10212

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 233
$ws.Range("I4").Value = 79.59999999999999
$ws.Range("K4").Value = 79.59999999999999
$ws.Range("M4").Value = 35.40000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1571.375
$ws.Range("I22").Value = 5150.5
$ws.Range("J22").Value = 378.33334
$ws.Range("K22").Value = 5150.5
$ws.Range("L22").Value = 378.33334
$ws.Range("M22").Value = -4800.5
$ws.Range("N22").Value = -1078.33334

$ws.Range("H132").Value = 1173
$ws.Range("I132").Value = 747.8125
$ws.Range("K132").Value = 2243.4375
$ws.Range("M132").Value = 286.5625

$ws.Range("H134").Value = 15152612
$ws.Range("I134").Value = 937.86206
$ws.Range("J134").Value = 125002250
$ws.Range("K134").Value = 2813.58618
$ws.Range("L134").Value = 375006750
$ws.Range("M134").Value = -278.5861800000002
$ws.Range("N134").Value = -375011820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1808.7142
$ws.Range("J5").Value = 2547.4443
$ws.Range("L5").Value = 7642.3329
$ws.Range("N5").Value = -7866.3329

$ws.Range("H80").Value = 7242.857
$ws.Range("I80").Value = 5900
$ws.Range("J80").Value = 7346.154
$ws.Range("K80").Value = 17700
$ws.Range("L80").Value = 22038.462
$ws.Range("M80").Value = -16764
$ws.Range("N80").Value = -23910.462

$ws.Range("H83").Value = 7242.857
$ws.Range("I83").Value = 5900
$ws.Range("J83").Value = 7346.154
$ws.Range("K83").Value = 53100
$ws.Range("L83").Value = 66115.386
$ws.Range("M83").Value = -48420
$ws.Range("N83").Value = -75475.386

$ws.Range("H122").Value = 1161.6904
$ws.Range("I122").Value = 363.85715
$ws.Range("J122").Value = 1321.2572
$ws.Range("K122").Value = 3274.71435
$ws.Range("L122").Value = 11891.3148
$ws.Range("M122").Value = -824.7143499999997
$ws.Range("N122").Value = -16791.3148

$ws.Range("H126").Value = 46097.082
$ws.Range("J126").Value = 4811
$ws.Range("L126").Value = 14433
$ws.Range("N126").Value = -24313

$ws.Range("H131").Value = 19309438
$ws.Range("I131").Value = 100202160
$ws.Range("K131").Value = 300606480
$ws.Range("M131").Value = -300601440

$ws.Range("H135").Value = 1808.7142
$ws.Range("J135").Value = 2547.4443
$ws.Range("L135").Value = 22926.9987
$ws.Range("N135").Value = -27996.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4089.5
$ws.Range("I70").Value = 3626.6667
$ws.Range("K70").Value = 3626.6667
$ws.Range("M70").Value = -3356.6667

$ws.Range("H73").Value = 4089.5
$ws.Range("I73").Value = 3626.6667
$ws.Range("K73").Value = 3626.6667
$ws.Range("M73").Value = -2690.6667

$ws.Range("H113").Value = 3048.7144
$ws.Range("I113").Value = 960.3333
$ws.Range("J113").Value = 4615
$ws.Range("K113").Value = 960.3333
$ws.Range("L113").Value = 4615
$ws.Range("M113").Value = 1209.6667
$ws.Range("N113").Value = -8955

$ws.Range("H132").Value = 3185.4443
$ws.Range("I132").Value = 2889.3333
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 8667.999899999999
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -6137.999899999999
$ws.Range("N132").Value = -19058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2677.7273
$ws.Range("I7").Value = 2745.5
$ws.Range("K7").Value = 2745.5
$ws.Range("M7").Value = -2633.5

$ws.Range("H22").Value = 1096.9584
$ws.Range("I22").Value = 304.57144
$ws.Range("J22").Value = 1423.2354
$ws.Range("K22").Value = 304.57144
$ws.Range("L22").Value = 1423.2354
$ws.Range("M22").Value = -9.571439999999996
$ws.Range("N22").Value = -2013.2354

$ws.Range("H27").Value = 1096.9584
$ws.Range("I27").Value = 304.57144
$ws.Range("J27").Value = 1423.2354
$ws.Range("K27").Value = 304.57144
$ws.Range("L27").Value = 1423.2354
$ws.Range("M27").Value = -197.57144
$ws.Range("N27").Value = -1637.2354

$ws.Range("H126").Value = 2677.7273
$ws.Range("I126").Value = 2745.5
$ws.Range("K126").Value = 8236.5
$ws.Range("M126").Value = -5766.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 649.125
$ws.Range("I107").Value = 556.0476
$ws.Range("K107").Value = 1668.1428
$ws.Range("M107").Value = 251.8571999999999

$ws.Range("H132").Value = 908.551
$ws.Range("I132").Value = 773.5278
$ws.Range("J132").Value = 1282.4615
$ws.Range("K132").Value = 2320.5834
$ws.Range("L132").Value = 3847.3845
$ws.Range("M132").Value = 209.4166
$ws.Range("N132").Value = -8907.3845
